$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44511
$ws.Range("M2").Value = 15
$ws.Range("N2").Value = 22000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 22000
$ws.Range("S2").Value = 1467

# Row 3
$ws.Range("D3").Value = 45079
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("S3").Value = 1000

# Row 4
$ws.Range("D4").Value = 44217
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = "`$/caja 18 kilos granel"
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44418
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("Q5").Value = "`$/caja 15 kilos granel"
$ws.Range("S5").Value = 533
$ws.Range("T5").Value = 15

# Row 6
$ws.Range("D6").Value = 45092
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 220
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("Q6").Value = "`$/caja 18 kilos granel"
$ws.Range("S6").Value = 889
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44208
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 24000
$ws.Range("P7").Value = 24000
$ws.Range("Q7").Value = "`$/caja 15 kilos granel"
$ws.Range("S7").Value = 1600
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 45085
$ws.Range("M8").Value = 110

# Row 9
$ws.Range("D9").Value = 45093
$ws.Range("M9").Value = 170
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15471
$ws.Range("S9").Value = 860

# Row 10
$ws.Range("D10").Value = 44411
$ws.Range("M10").Value = 210
$ws.Range("N10").Value = 8000
$ws.Range("O10").Value = 8000
$ws.Range("P10").Value = 8000
$ws.Range("Q10").Value = "`$/bandeja 8 kilos"
$ws.Range("S10").Value = 1000
$ws.Range("T10").Value = 8

# Row 11
$ws.Range("D11").Value = 44966
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 4
$ws.Range("N11").Value = 250000
$ws.Range("O11").Value = 250000
$ws.Range("P11").Value = 250000
$ws.Range("Q11").Value = "`$/bins (400 kilos)"
$ws.Range("S11").Value = 625
$ws.Range("T11").Value = 400

# Row 12
$ws.Range("D12").Value = 44966
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("S12").Value = 833

# Row 13
$ws.Range("D13").Value = 45086
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 16000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 16000
$ws.Range("S13").Value = 889

# Row 15
$ws.Range("D15").Value = 44427
$ws.Range("M15").Value = 55
$ws.Range("N15").Value = 7000
$ws.Range("O15").Value = 7000
$ws.Range("P15").Value = 7000
$ws.Range("Q15").Value = "`$/caja 15 kilos granel"
$ws.Range("S15").Value = 467
$ws.Range("T15").Value = 15

# Row 16
$ws.Range("D16").Value = 44264
$ws.Range("L16").Value = "Calibre 100"
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 20000
$ws.Range("Q16").Value = "`$/caja 18 kilos embalada"
$ws.Range("S16").Value = 1111

# Row 17
$ws.Range("D17").Value = 45090
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 140
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("Q17").Value = "`$/caja 18 kilos granel"
$ws.Range("S17").Value = 889

# Row 18
$ws.Range("D18").Value = 45089
$ws.Range("M18").Value = 100

# Row 19
$ws.Range("D19").Value = 44601
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 28000
$ws.Range("O19").Value = 28000
$ws.Range("P19").Value = 28000
$ws.Range("Q19").Value = "`$/caja 18 kilos granel"
$ws.Range("S19").Value = 1556
$ws.Range("T19").Value = 18

# Row 20
$ws.Range("D20").Value = 44392
$ws.Range("L20").Value = "Especial"
$ws.Range("M20").Value = 500
$ws.Range("N20").Value = 7000
$ws.Range("O20").Value = 7000
$ws.Range("P20").Value = 7000
$ws.Range("Q20").Value = "`$/bandeja 8 kilos"
$ws.Range("S20").Value = 875
$ws.Range("T20").Value = 8
